$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.1854658009979175

# Row 3
$ws.Range("D3").Value = 6.37490644099671

# Row 4
$ws.Range("D4").Value = 10.85041264000029

# Row 5
$ws.Range("D5").Value = 25.4916347790022

# Row 6
$ws.Range("D6").Value = 64.11864413600051

# Row 7
$ws.Range("D7").Value = 0.7694899710004393

# Row 8
$ws.Range("D8").Value = 13.79679634100103

# Row 9
$ws.Range("D9").Value = 27.5478598539994

# Row 10
$ws.Range("D10").Value = 73.90364072800003

# Row 11
$ws.Range("D11").Value = 193.8681986440024

# Row 12
$ws.Range("D12").Value = 23.33660387700002
$ws.Range("E12").Value = 13

# Row 13
$ws.Range("D13").Value = 38.41004757100018
$ws.Range("E13").Value = 11

# Row 14
$ws.Range("D14").Value = 8.519867435999913
$ws.Range("E14").Value = 11

# Row 15
$ws.Range("B15").Value = 321
$ws.Range("C15").Value = 1634
$ws.Range("D15").Value = 0.01353961700078798

# Row 16
$ws.Range("B16").Value = 808
$ws.Range("C16").Value = 7842
$ws.Range("D16").Value = 0.100543580003432

# Row 17
$ws.Range("D17").Value = 0.152428344998043

# Row 18
$ws.Range("D18").Value = 0.00380709699675208

# Row 19
$ws.Range("D19").Value = 0.02201738800067687

# Row 20
$ws.Range("B20").Value = 703
$ws.Range("C20").Value = 6148
$ws.Range("D20").Value = 0.1145449559990084

# Row 21
$ws.Range("D21").Value = 0.006563068003742956

# Row 22
$ws.Range("D22").Value = 0.02292848400247749

# Row 23
$ws.Range("D23").Value = 0.1595110019989079

# Row 24
$ws.Range("B24").Value = 718
$ws.Range("C24").Value = 7453
$ws.Range("D24").Value = 0.1674030559952371

# Row 25
$ws.Range("D25").Value = 0.1182485919998726

# Row 26
$ws.Range("D26").Value = 0.3190244579964201

# Row 27
$ws.Range("D27").Value = 0.0305055179996998

# Row 28
$ws.Range("D28").Value = 0.339555175996793

# Row 29
$ws.Range("D29").Value = 7.186975135999091

# Row 30
$ws.Range("D30").Value = 6.565853181004059

# Row 31
$ws.Range("D31").Value = 20.78300001799653

# Row 32
$ws.Range("D32").Value = 8.217602974000329

# Row 33
$ws.Range("D33").Value = 440.8838805600026
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = "COMPLETE"

# Row 34
$ws.Range("D34").Value = 56.34603770499962

# Row 35
$ws.Range("D35").Value = 0.418952038999123

# Row 36
$ws.Range("D36").Value = 23.83931398599816

# Row 37
$ws.Range("D37").Value = 0.2540792230000193

# Row 38
$ws.Range("D38").Value = 2.287235221003357

# Row 39
$ws.Range("D39").Value = 2.423461952996149

# Row 40
$ws.Range("D40").Value = 0.1630570299967076

# Row 41
$ws.Range("D41").Value = 25.76155824199668

# Row 42
$ws.Range("D42").Value = 70.51578570099809

# Row 43
$ws.Range("D43").Value = 26.26395418599714

# Row 44
$ws.Range("D44").Value = 83.53621197000029

# Row 45
$ws.Range("D45").Value = 14.56240587699722

# Row 46
$ws.Range("D46").Value = 3.59251386699907

# Row 47
$ws.Range("D47").Value = 49.79581049199624

# Row 48
$ws.Range("D48").Value = 65.47249596600159

# Row 49
$ws.Range("D49").Value = 267.5666719950023
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = "COMPLETE"

# Row 50
$ws.Range("B50").Value = 1618
$ws.Range("C50").Value = 17262
$ws.Range("D50").Value = 169.7970782269986

# Row 51
$ws.Range("D51").Value = 94.84138439600065
